# Bulk-upload template: the "leadowneremail" column is no longer needed.
# Delete the entire column A, which shifts source/firstname/lastname/email/contact
# (and the stray " " helper value in column F) one column to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireColumn.Delete()

# Leave the active selection where Excel would land after the delete/edit session.
$ws.Range("D7").Select()
